# "changed icon names and updated spreadsheet"
# Adds a new "Icon Src" column (I) to Sheet1 with per-company icon paths,
# matching each existing data row by its client name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, bold like the rest of row 1.
$ws.Range("I1").Value = "Icon Src"
$ws.Range("I1").Font.Bold = $true

# One icon path per company row (rows without a listed icon keep column I empty).
$ws.Range("I2").Value  = "/assets/icons/SARU_logo.svg"
$ws.Range("I6").Value  = "/assets/icons/afmetco-logo.png"
$ws.Range("I7").Value  = "/assets/icons/sagoodnews-pernod-ricard.png"
$ws.Range("I8").Value  = "/assets/icons/Nestle-Logo.png"
$ws.Range("I11").Value = "/assets/icons/Country-Bird-Holdings.png"
$ws.Range("I12").Value = "/assets/icons/karan_beef.png"
$ws.Range("I13").Value = "/assets/icons/afgri-logo-for-posts.jpg"
$ws.Range("I14").Value = "/assets/icons/homechoice.co.za.png"
$ws.Range("I15").Value = "/assets/icons/telkom-logo.jpg"
$ws.Range("I16").Value = "/assets/icons/vodacom.jpeg"
$ws.Range("I17").Value = "/assets/icons/woolworths.jpeg"
$ws.Range("I18").Value = "/assets/icons/mrprice.jpeg"
$ws.Range("I19").Value = "/assets/icons/shoprite_checkers_logo.jpeg"

# Page setup for the (now wider) sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Zoom the view out a bit and park the selection just past the used range.
$excel.ActiveWindow.Zoom = 85
$ws.Range("I20").Select() | Out-Null
